$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.420.17'
$ws.Range('E2').Value = '  -2.71%  '
$ws.Range('D3').Value = '2.895.27'
$ws.Range('E3').Value = '  -3.81%  '
$ws.Range('D5').Value = '''585.88'
$ws.Range('E5').Value = '  -1.36%  '
$ws.Range('D6').Value = '''146.63'
$ws.Range('E6').Value = '  -0.33%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('E8').Value = '  -2.65%  '
$ws.Range('D9').Value = '2.894.03'
$ws.Range('E9').Value = '  -3.85%  '
$ws.Range('D10').Value = '''6.64'
$ws.Range('E10').Value = '  +6.31%  '
$ws.Range('E11').Value = '  -3.22%  '
$ws.Range('D12').Value = '''0.447'
$ws.Range('E12').Value = '  -2.21%  '
$ws.Range('D13').Value = '''0.0000224'
$ws.Range('E13').Value = '  -3.42%  '
$ws.Range('D14').Value = '''34.16'
$ws.Range('E14').Value = '  -0.65%  '
$ws.Range('E15').Value = '  +0.31%  '
$ws.Range('D16').Value = '3.375.30'
$ws.Range('E16').Value = '  -5.98%  '
$ws.Range('D17').Value = '''6.80'
$ws.Range('E17').Value = '  -2.43%  '
$ws.Range('D18').Value = '60.436.43'
$ws.Range('E18').Value = '  -2.62%  '
$ws.Range('D19').Value = '2.893.78'
$ws.Range('E19').Value = '  -3.61%  '
$ws.Range('D20').Value = '''424.74'
$ws.Range('E20').Value = '  -4.78%  '
$ws.Range('D21').Value = '''13.60'
$ws.Range('E21').Value = '  -4.04%  '
$ws.Range('D22').Value = '''0.670'
$ws.Range('E22').Value = '  -2.53%  '
$ws.Range('D23').Value = '''7.11'
$ws.Range('E23').Value = '  -3.86%  '
$ws.Range('D24').Value = '''80.88'
$ws.Range('E24').Value = '  -1.67%  '
$ws.Range('D25').Value = '''11.00'
$ws.Range('E25').Value = '  -0.15%  '
$ws.Range('E26').Value = '  -3.37%  '
$ws.Range('D27').Value = '''11.77'
$ws.Range('E27').Value = '  -2.55%  '
$ws.Range('E28').Value = '  +0.01%  '
$ws.Range('B29').Value = 'NEARProtocol'
$ws.Range('C29').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D29').Value = '''7.28'
$ws.Range('E29').Value = '  +1.05%  '
$ws.Range('B30').Value = 'FirstDigitalUSD'
$ws.Range('C30').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D30').Value = '''0.999'
$ws.Range('E30').Value = '  -0.09%  '
$ws.Range('D31').Value = '''2.18'
$ws.Range('E31').Value = '  +3.79%  '
$ws.Range('D32').Value = '''2.62'
$ws.Range('E32').Value = '  -3.17%  '
$ws.Range('D33').Value = '''26.51'
$ws.Range('E33').Value = '  -3.32%  '
$ws.Range('E34').Value = '  -3.95%  '
$ws.Range('D35').Value = '0.0₃0834'
$ws.Range('E35').Value = '  -1.98%  '
$ws.Range('E36').Value = '  -2.28%  '
$ws.Range('D37').Value = '''5.66'
$ws.Range('E37').Value = '  -2.84%  '
$ws.Range('D38').Value = '''49.38'
$ws.Range('E38').Value = '  -1.71%  '
$ws.Range('E39').Value = '  -1.05%  '
$ws.Range('D40').Value = '''2.93'
$ws.Range('E40').Value = '  -1.41%  '
$ws.Range('E41').Value = '  -1.35%  '
$ws.Range('D42').Value = '''8.75'
$ws.Range('E42').Value = '  -3.05%  '
$ws.Range('E43').Value = '  +0.93%  '
$ws.Range('D44').Value = '''41.14'
$ws.Range('E44').Value = '  -1.11%  '
$ws.Range('D45').Value = '''0.0345'
$ws.Range('E45').Value = '  -1.77%  '
$ws.Range('D46').Value = '''371.51'
$ws.Range('E46').Value = '  -5.65%  '
$ws.Range('B47').Value = 'Maker'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D47').Value = '2.649.04'
$ws.Range('E47').Value = '  -2.79%  '
$ws.Range('B48').Value = 'Monero'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D48').Value = '''132.94'
$ws.Range('E48').Value = '  -1.07%  '
$ws.Range('E49').Value = '  +0.00%  '
$ws.Range('D50').Value = '''25.14'
$ws.Range('E50').Value = '  +5.31%  '
$ws.Range('E51').Value = '  -1.00%  '
